# fix(DOCX): Encode table captions; update template
#
# 1. TableCaption style: add spacing (before=28.35pt/567tw, after=14.15pt/283tw),
#    left alignment, and a 12pt (sz=24) run size.
# 2. ImageCaption style: add a 12pt (sz=24) run size.
# 3. New FigureCaption style (based on Figure) with spacing
#    before=14.15pt/283tw, after=28.35pt/567tw.

$d = $word.ActiveDocument
$styles = $d.Styles

# --- TableCaption -----------------------------------------------------
$tableCaption = $styles.Item("TableCaption")
$tableCaption.ParagraphFormat.SpaceBefore = 28.35
$tableCaption.ParagraphFormat.SpaceAfter = 14.15
$tableCaption.ParagraphFormat.Alignment = 0
$tableCaption.Font.Size = 12

# --- ImageCaption -------------------------------------------------------
$imageCaption = $styles.Item("ImageCaption")
$imageCaption.Font.Size = 12

# --- FigureCaption (new style) ------------------------------------------
$figureCaption = $styles.Add("Figure Caption", 1)
$figureCaption.QuickStyle = $true
$figureCaption.BaseStyle = $styles.Item("Figure")
$figureCaption.ParagraphFormat.SpaceBefore = 14.15
$figureCaption.ParagraphFormat.SpaceAfter = 28.35

Write-Host "Styles updated: TableCaption, ImageCaption, FigureCaption(new)"
